$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Columns("D:E").Insert()

# Copy cell formatting (number formats / styles) from column F (the old column D, now
# shifted right by the insert) into the two newly inserted columns D:E, block by block
# (the three financial-statement sections: Income Statement, Balance Sheet, Cash Flow).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new quarter columns (D = 2018-12-31, E = 2018-09-30) with data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 146100
$ws.Range("E8").Value = 148200
$ws.Range("D9").Value = 116700
$ws.Range("E9").Value = 116700
$ws.Range("D10").Value = 29400
$ws.Range("E10").Value = 31500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 132500
$ws.Range("E17").Value = 131800
$ws.Range("D18").Value = 13600
$ws.Range("E18").Value = 16400
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 19400
$ws.Range("E21").Value = 22000
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 13600
$ws.Range("E23").Value = 16400
$ws.Range("D24").Value = 3300
$ws.Range("E24").Value = 3600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 10300
$ws.Range("E26").Value = 12800
$ws.Range("D27").Value = 10300
$ws.Range("E27").Value = 12800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 10300
$ws.Range("E33").Value = 12800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 10300
$ws.Range("E35").Value = 12800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 96500
$ws.Range("E41").Value = 90100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 62400
$ws.Range("E43").Value = 65100
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 51100
$ws.Range("E45").Value = 48900
$ws.Range("D46").Value = 210000
$ws.Range("E46").Value = 204100
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 62900
$ws.Range("E48").Value = 64800
$ws.Range("D49").Value = 39900
$ws.Range("E49").Value = 40000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1000
$ws.Range("E52").Value = 1000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 313700
$ws.Range("E54").Value = 309900
$ws.Range("D57").Value = 9900
$ws.Range("E57").Value = 9100
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 107500
$ws.Range("E59").Value = 106500
$ws.Range("D60").Value = 117400
$ws.Range("E60").Value = 115600
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 4100
$ws.Range("E62").Value = 4300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 121500
$ws.Range("E66").Value = 119900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 493300
$ws.Range("E72").Value = 483000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 192300
$ws.Range("E76").Value = 190100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 10300
$ws.Range("E81").Value = 12800
$ws.Range("D83").Value = 5800
$ws.Range("E83").Value = 5600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 19300
$ws.Range("E89").Value = 25800
$ws.Range("D91").Value = -3800
$ws.Range("E91").Value = -3500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -3800
$ws.Range("E94").Value = -3500
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -9100
$ws.Range("E100").Value = -5500
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 6400
$ws.Range("E102").Value = 16800
